# Apply edits to parallel.xlsx worksheet: extend table from columns A:O to A:Q
# by adding two new columns (P, Q) with header values 14 and 15 (matching the
# header row style), filling the new data-row cells with the value 2, and
# swapping the I/K and M/O column values for all data rows (2-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1 = 14, Q1 = 15, matching the style of the
# existing header cells (e.g. O1). ---
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("P1").Value2 = 14
$ws.Range("Q1").PasteSpecial(-4122)
$ws.Range("Q1").Value2 = 15

# --- Data rows (2-25): swap I<->K and M<->O values, then fill P and Q with 2 ---
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I
    $kVal = $ws.Cells.Item($r, 11).Value2  # column K
    $mVal = $ws.Cells.Item($r, 13).Value2  # column M
    $oVal = $ws.Cells.Item($r, 15).Value2  # column O

    $ws.Cells.Item($r, 9).Value2 = $kVal   # I = old K
    $ws.Cells.Item($r, 11).Value2 = $iVal  # K = old I
    $ws.Cells.Item($r, 13).Value2 = $oVal  # M = old O
    $ws.Cells.Item($r, 15).Value2 = $mVal  # O = old M

    $ws.Cells.Item($r, 16).Value2 = 2      # column P
    $ws.Cells.Item($r, 17).Value2 = 2      # column Q
}
